# Added empty cols test
# Select column D and insert a new (blank) column before it, shifting the
# existing D and E columns one place to the right (D->E, E->F). This mirrors
# a user selecting the "D" column header and choosing Insert > Sheet Columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column C's width so the freshly inserted column D can inherit it,
# the same way Excel copies formatting from the column to the left.
$colCWidth = $ws.Columns.Item(3).ColumnWidth

# Remember the shape's target cell-relative position before the sheet
# layout changes, so it can be re-anchored one column to the right below.
$shape = $ws.Shapes.Item(1)

# Perform the actual column insert.
$ws.Columns("D:D").Insert()

# New column D should look like old column C (format-wise / width-wise).
$ws.Columns.Item(4).ColumnWidth = $colCWidth

# Defined names that pointed into column E now need to point at column F
# (everything from column D onward shifted right by one column).
foreach ($n in $wb.Names) {
    $ref = $n.RefersTo
    if ($ref.Contains("`$E`$")) {
        $n.RefersTo = $ref.Replace("`$E`$", "`$F`$")
    }
}

# Cell comments don't automatically follow the column shift, so move them
# from the old D1/D10 positions to the new E1/E10 positions.
$comment1 = $ws.Range("D1").Comment
if ($comment1 -ne $null) {
    $text1 = $comment1.Text()
    $comment1.Delete()
    $ws.Range("E1").AddComment($text1)
}

$comment10 = $ws.Range("D10").Comment
if ($comment10 -ne $null) {
    $text10 = $comment10.Text()
    $comment10.Delete()
    $ws.Range("E10").AddComment($text10)
}

# The floating text box also doesn't move with the column insert, so
# re-anchor it one column to the right (same row), preserving its
# position relative to the grid.
$targetCell = $ws.Cells.Item(11, 7)
$shape.Left = $targetCell.Left
$shape.Top = $targetCell.Top

# Leave the whole new column D selected, as if the user had just inserted it
# via the column header.
$ws.Range("D1:D1048576").Select()
